# Finished draft of introduction to physarum.
# Move the "PP applicatons" block (old rows 26-32) down to rows 37-43 to make
# room, relocate the "Introduction to PP" block (old rows 9-14) down to rows
# 27-32, and append a new "motivation for the thesis" note at row 45.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Introduction")
$ws.Activate()

# Hyperlink ranges do not automatically follow cells that get cut/pasted, so
# drop the stale collection up front and rebuild it once every cell is in
# its final resting place.
$ws.Hyperlinks.Delete()

# 1) Slide the "PP applicatons" bibliography block from rows 26-32 down to
#    rows 37-43 (an 11-row move) to make space for the relocated
#    "Introduction to PP" block.
$ws.Range("E26:M32").Cut($ws.Range("E37"))

# 2) Move the "Introduction to PP" / taxonomy block from rows 9-14 down to
#    rows 27-32 (an 18-row move). Row 9 keeps its "PP and Natural
#    computing" heading in column C.
$ws.Range("E9:J14").Cut($ws.Range("E27"))

# 3) New note added below the restructured outline.
$ws.Range("E45").Value = "motivation for the thesis"

# 4) Re-create the 10 hyperlinks at their new homes, preserving the original
#    target URLs.
$ws.Hyperlinks.Add($ws.Range("M20"), "http://www.sciencedirect.com/science/article/pii/S1672652911600164")
$ws.Hyperlinks.Add($ws.Range("M21"), "http://science.sciencemag.org/content/327/5964/439.short")
$ws.Hyperlinks.Add($ws.Range("M22"), "http://www.sciencedirect.com/science/article/pii/S0303264711000803")
$ws.Hyperlinks.Add($ws.Range("M23"), "http://www.sciencedirect.com/science/article/pii/S1878778911000305")
$ws.Hyperlinks.Add($ws.Range("M24"), "http://www.sciencedirect.com/science/article/pii/S0303264706001687")
$ws.Hyperlinks.Add($ws.Range("M38"), "http://www.hindawi.com/journals/tswj/2014/487069/abs/")
$ws.Hyperlinks.Add($ws.Range("M40"), "http://search.ebscohost.com/login.aspx?direct=true&profile=ehost&scope=site&authtype=crawler&jrnl=15487199&AN=92711421&h=XVLV%2FphLLXPzP%2Fck30zliftFLCgxX%2F5FjC%2FH70GX70vdRq8OmVX3dzxb8D1lIGkjXWJeeE%2B8gwBCcSkRw9TINA%3D%3D&crl=c")
$ws.Hyperlinks.Add($ws.Range("M41"), "http://ieeexplore.ieee.org/xpls/abs_all.jsp?arnumber=6684158")
$ws.Hyperlinks.Add($ws.Range("M42"), "http://dl.acm.org/citation.cfm?id=2744528")
$ws.Hyperlinks.Add($ws.Range("M43"), "http://journals.plos.org/plosone/article?id=10.1371/journal.pone.0066732")

# Adding hyperlinks re-applies Excel's built-in "Hyperlink" cell style; put
# the original centered style back on every cell that carries one so the
# bibliography list keeps its original look.
$ws.Range("M25").Copy()
$ws.Range("M20:M24").PasteSpecial(-4122)
$ws.Range("M38").PasteSpecial(-4122)
$ws.Range("M40:M43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Park the selection where the author's edit left it.
$ws.Range("C45").Select()
